# Loan RBI, Variable Instalments
#
# Inserts a new (blank) "Variable Instalment" column into the
# "Repayment Schedule" sheet between the existing "In Advance" (M) and
# "Late" (N) columns, shifting "Late"/"heading"/"Outstanding" one column
# to the right (N->O, O->P, P->Q), and leaves the Repayment Schedule tab
# as the active/selected sheet & cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a brand-new blank column at N (14th column); everything from N
# onward (N, O, P) shifts right to (O, P, Q), carrying its data/styles.
$ws.Columns.Item(14).Insert()

# Match the original column's display width (Excel "characters" units map
# to the stored XML width with a +5/6 padding offset, so 9.1666... in
# ColumnWidth terms serializes back out as width="10").
$ws.Columns.Item(14).ColumnWidth = 9.166666666666666

# Make "Repayment Schedule" the active sheet/tab and select T11, matching
# where the editor was working after inserting the column.
$ws.Activate()
$ws.Range("T11").Select()
